$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Parent company" column (B) and the "Location County/City"
# column (originally E, now D after the first delete) - additional
# companies were sent for questionaire, these two columns are no longer
# part of the tracked sheet.
$ws.Range("B:B").EntireColumn.Delete()
$ws.Range("D:D").EntireColumn.Delete()

$ws.Range("A2:I2").Select()
